$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Heading: "Some Class Name" -> "RuleSubView"
# ---------------------------------------------------------------------
$headingRange = $d.Paragraphs.Item(1).Range
$headingRange.Find.Execute("Some Class Name", $true, $false, $false, $false, $false, `
    $true, 1, $false, "RuleSubView", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. CRC table, first row / cell: "Class Name" -> "RuleSubView"
#    (append a sentinel "X" so we can later carve out a clean, collapsed
#    insertion point for the relocated _GoBack bookmark without landing
#    exactly on the paragraph-end mark)
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$cellRange = $table.Cell(1, 1).Range
$cellRange.Find.Execute("Class Name", $true, $false, $false, $false, $false, `
    $true, 1, $false, "RuleSubViewX", 2) | Out-Null

$cellRange = $table.Cell(1, 1).Range
$newNameLength = "RuleSubView".Length
$bookmarkPos = $cellRange.Start + $newNameLength

# ---------------------------------------------------------------------
# 3. Relocate the reserved "_GoBack" bookmark from its old spot (inside
#    the "Some paragraph about this class" run-pair under Description)
#    to right after the new "RuleSubView" text in the CRC table.
# ---------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$goBackRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Remove the sentinel "X" now that the bookmark is anchored just before it.
$sentinel = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinel.Delete()
